$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 31; $r -le 74; $r++) {
    $ws.Range("I$r").Value = 63.85925373134329
}
